$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "KKE - Communist Party of Greece  (Kommounistiko Komma Elladas, KKE)"
$ws.Range("C1").Value = "ND - New Democracy  (Nea Dimokratia, ND)"
$ws.Range("D1").Value = "PASOK - Pan-Hellenic Socialist Movement  (Panellinio Sosialistiko Kinima, PASOK)"
$ws.Range("E1").Value = "POLAN - Political Spring  (Politiki Anixi, POLAN)"
$ws.Range("F1").Value = "DIKKI - Democratic Social Movement/ Dimokratiko Koinoniko Kinima (DIKKI)  (Dimokratiko Koinoniko Kinima, DIKKI)"
$ws.Range("G1").Value = "SYN - Coalition of the Left and Progress  (Synaspismos tis Aristeras kai tis Proodou, SYN)"
$ws.Range("H1").Value = "SYRIZA - Coalition of the Radical Left  (Synaspismos tis Rizospastikis Aristeras, SYRIZA)"
$ws.Range("I1").Value = "LAOS - Popular Orthodox Rally (Laikos Orthodoxos Synagermos, LAOS)"
$ws.Range("J1").Value = "EL, GS - Greek Solution (Elliniki Lysi, EL, GS)"
$ws.Range("K1").Value = "KINAL - Movement for Change"
$ws.Range("L1").Value = "MeRA25 - European Realistic Disobedience Front"
$ws.Range("M1").Value = "ANEL - Independent Hellenes (Anexartitoi Ellines, ANEL)"
$ws.Range("N1").Value = "DIMAR - Democratic Left  (Dimokratiki Aristera, DIMAR)"
$ws.Range("O1").Value = "HA - Golden Dawn (Chrysi Avgi, HA)"
$ws.Range("P1").Value = "EK - Union of Centrists (Enosi Kentroon, EK)"
$ws.Range("Q1").Value = "TP - The River (To Potami , TP)"
